# Applies ARIMA bug-fix value updates to the "y_value" (column B) series
# in the "y_fitted_on_begin_2016" and "y_fitted_on_begin_2021" worksheets.

$wb = $excel.ActiveWorkbook

$sheet1Updates = @(
    @{Row=2; Value=0.3872817209427971},
    @{Row=3; Value=34.19454045665108},
    @{Row=4; Value=34.25890035670417},
    @{Row=5; Value=35.28411899868561},
    @{Row=6; Value=34.90888235148348},
    @{Row=7; Value=35.8080978059046},
    @{Row=8; Value=36.50228149299173},
    @{Row=9; Value=37.37574955658741},
    @{Row=10; Value=37.83956121323219},
    @{Row=11; Value=37.18372722637498},
    @{Row=12; Value=37.61667970049238},
    @{Row=13; Value=37.43206691559087},
    @{Row=14; Value=38.17297174466074},
    @{Row=15; Value=38.29984990300353},
    @{Row=16; Value=38.87796663796539},
    @{Row=17; Value=37.88358406077649},
    @{Row=18; Value=38.23694983927675},
    @{Row=19; Value=38.09661119454078},
    @{Row=20; Value=39.46683324454744},
    @{Row=21; Value=40.51985044788169},
    @{Row=22; Value=42.35301950779714},
    @{Row=23; Value=43.12020836389738},
    @{Row=24; Value=44.19624298393092},
    @{Row=25; Value=45.11610927507318},
    @{Row=26; Value=45.09192737424289},
    @{Row=27; Value=44.50214069318987},
    @{Row=28; Value=44.83021464369263},
    @{Row=29; Value=44.43051743243277},
    @{Row=30; Value=44.39543964485716},
    @{Row=31; Value=46.26983162794905},
    @{Row=32; Value=46.37216453285728},
    @{Row=33; Value=45.42103128356187},
    @{Row=34; Value=45.82942697170549},
    @{Row=35; Value=46.68842901902573},
    @{Row=36; Value=47.08599649031021},
    @{Row=37; Value=48.20969547682593},
    @{Row=38; Value=48.17102984136208}
)

$sheet3Updates = @(
    @{Row=2; Value=0.3386138405455644},
    @{Row=3; Value=34.14587257625384},
    @{Row=4; Value=34.21023247630693},
    @{Row=5; Value=35.23545111828837},
    @{Row=6; Value=34.86021447108624},
    @{Row=7; Value=35.75942992550737},
    @{Row=8; Value=36.45361361259449},
    @{Row=9; Value=37.32708167619018},
    @{Row=10; Value=37.79089333283495},
    @{Row=11; Value=37.13505934597775},
    @{Row=12; Value=37.56801182009515},
    @{Row=13; Value=37.38339903519363},
    @{Row=14; Value=38.12430386426351},
    @{Row=15; Value=38.2511820226063},
    @{Row=16; Value=38.82929875756815},
    @{Row=17; Value=37.83491618037925},
    @{Row=18; Value=38.18828195887952},
    @{Row=19; Value=38.04794331414355},
    @{Row=20; Value=39.41816536415021},
    @{Row=21; Value=40.47118256748445},
    @{Row=22; Value=42.30435162739991},
    @{Row=23; Value=43.07154048350014},
    @{Row=24; Value=44.14757510353369},
    @{Row=25; Value=45.06744139467595},
    @{Row=26; Value=45.04325949384566},
    @{Row=27; Value=44.45347281279264},
    @{Row=28; Value=44.7815467632954},
    @{Row=29; Value=44.38184955203553},
    @{Row=30; Value=44.34677176445993},
    @{Row=31; Value=46.22116374755181},
    @{Row=32; Value=46.32349665246005},
    @{Row=33; Value=45.37236340316463},
    @{Row=34; Value=45.78075909130825},
    @{Row=35; Value=46.63976113862849},
    @{Row=36; Value=47.03732860991298},
    @{Row=37; Value=48.1610275964287},
    @{Row=38; Value=48.12236196096485},
    @{Row=39; Value=48.08801513052038},
    @{Row=40; Value=48.0095038089868},
    @{Row=41; Value=48.77018015601529},
    @{Row=42; Value=50.44628275787453},
    @{Row=43; Value=48.77547061304323}
)

$wsFitted2016 = $wb.Worksheets.Item("y_fitted_on_begin_2016")
foreach ($u in $sheet1Updates) {
    $wsFitted2016.Cells.Item($u.Row, 2).Value = $u.Value
}

$wsFitted2021 = $wb.Worksheets.Item("y_fitted_on_begin_2021")
foreach ($u in $sheet3Updates) {
    $wsFitted2021.Cells.Item($u.Row, 2).Value = $u.Value
}
